$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2:64) down by one row (3:65), reading each source
# value explicitly (bottom-up so sources aren't clobbered before they're read).
for ($r = 64; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($dest, 4).Value = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($dest, 5).Value = $ws.Cells.Item($r, 5).Value()
}

# New row 2: a July (Mes=7) entry for "Dia" 1
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 8952.83
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 2025
$ws.Range("E2").Value = "07/2025"

# Corrections to a handful of shifted rows' totals
$ws.Range("B15").Value = 14906.35
$ws.Range("B19").Value = 46214.09
$ws.Range("B23").Value = 9530.559999999999
$ws.Range("B24").Value = 104974.48
